$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new rows of data
$ws.Range("A3").Value = "Heather"
$ws.Range("B3").Value = "unplanned;POM25;POM26"
$ws.Range("C3").Value = "5;25;10"

$ws.Range("A4").Value = "Monica"
$ws.Range("B4").Value = "Admin;24BES;24Spend plan"
$ws.Range("C4").Value = "15;8;17"

# Update the selection to a multi-area selection matching the diff:
# the original A3:C4 block plus the newly touched cell F15.
$excel.Union($ws.Range("A3:C4"), $ws.Range("F15")).Select()
